# soulbase.xlsx edit: add first class soul
# - Updates the "race" comment on I1 to insert a new "野兽" (beast) race
#   option and renumber the old 矮人/精灵/兽人 entries to 人类/兽人/精灵.
# - Renumbers the first two soul rows (id 1/2 -> 1000/1001), tweaks their
#   stat columns, and appends 13 brand-new soul rows (1002, 2000-2002,
#   3000-3002, 4000-4002, 5000-5002) each with its own avatar/name/desc.
# - Widens column C and updates the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the race comment on I1.
# ---------------------------------------------------------------------
$nl = [char]10
$raceComment = "王朋:" + $nl + "种族" + $nl + $nl + "野兽1" + $nl + "人类2" + $nl + "兽人3" + $nl + "精灵4" + $nl + "亡灵5" + $nl + "龙族6" + $nl + "恶魔7" + $nl + "神族8" + $nl + $nl
$ws.Range("I1").Comment.Text($raceComment)

# ---------------------------------------------------------------------
# 2. Full target data for rows 2-16 (row 1 is the header, untouched).
#    Columns: A id, B icon, C canmutate, D avatar, E profession, F name,
#    G desc, H copper, I race, J rarityclass, K elemtype, L skillid,
#    M MaxHP, N Attack, O Defense, P Heal, Q Speed, R Dodge, S Crit,
#    T fatherrace, U fatherid, V motherrace, W motherid.
# ---------------------------------------------------------------------
$data = @(
  @(1000, "soul.png", 1, "sk_human,charactor/face.png,charactor/head.png,charactor/body.png,charactor/arm.png,charactor/hand.png,charactor/arm.png,charactor/hand.png,charactor/leg.png,charactor/leg.png,", 0, "SOUL_NAME_1000", "SOUL_DESC_1000", 1, 1, 1, 0, 3000, 50, 30, 10, 10, 10, 0, 0, 1, 0, 1, 0),
  @(1001, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor2/leg.png,", 0, "SOUL_NAME_1001", "SOUL_DESC_1001", 1, 1, 1, 0, 3000, 50, 40, 0, 10, 10, 0, 0, 2, 0, 4, 0),
  @(1002, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor3/leg.png,", 0, "SOUL_NAME_1002", "SOUL_DESC_1002", 1, 1, 1, 0, 3000, 50, 20, 20, 10, 10, 0, 0, 2, 0, 5, 0),
  @(2000, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor4/leg.png,", 0, "SOUL_NAME_2000", "SOUL_DESC_2000", 1, 2, 1, 0, 3000, 50, 40, 0, 10, 10, 0, 0, 2, 0, 2, 0),
  @(2001, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor5/leg.png,", 0, "SOUL_NAME_2001", "SOUL_DESC_2001", 1, 2, 1, 0, 3000, 50, 30, 10, 10, 10, 0, 0, 1, 0, 5, 0),
  @(2002, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor6/leg.png,", 0, "SOUL_NAME_2002", "SOUL_DESC_2002", 1, 2, 1, 0, 3000, 50, 40, 0, 10, 10, 0, 0, 4, 0, 3, 0),
  @(3000, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor7/leg.png,", 0, "SOUL_NAME_3000", "SOUL_DESC_3000", 1, 3, 1, 0, 3000, 50, 30, 10, 10, 10, 0, 0, 3, 0, 3, 0),
  @(3001, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor8/leg.png,", 0, "SOUL_NAME_3001", "SOUL_DESC_3001", 1, 3, 1, 0, 3000, 50, 20, 20, 10, 10, 0, 0, 1, 0, 2, 0),
  @(3002, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor9/leg.png,", 0, "SOUL_NAME_3002", "SOUL_DESC_3002", 1, 3, 1, 0, 3000, 50, 30, 10, 20, 10, 0, 0, 4, 0, 5, 0),
  @(4000, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor10/leg.png,", 0, "SOUL_NAME_4000", "SOUL_DESC_4000", 1, 4, 1, 0, 3000, 50, 30, 10, 10, 10, 0, 0, 4, 0, 4, 0),
  @(4001, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor11/leg.png,", 0, "SOUL_NAME_4001", "SOUL_DESC_4001", 1, 4, 1, 0, 3000, 50, 10, 30, 10, 10, 0, 0, 1, 0, 3, 0),
  @(4002, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor12/leg.png,", 0, "SOUL_NAME_4002", "SOUL_DESC_4002", 1, 4, 1, 0, 3000, 50, 40, 0, 10, 10, 0, 0, 2, 0, 5, 0),
  @(5000, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor13/leg.png,", 0, "SOUL_NAME_5000", "SOUL_DESC_5000", 1, 5, 1, 0, 3000, 50, 40, 10, 10, 10, 0, 0, 5, 0, 5, 0),
  @(5001, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor14/leg.png,", 0, "SOUL_NAME_5001", "SOUL_DESC_5001", 1, 5, 1, 0, 3000, 50, 40, 10, 10, 10, 0, 0, 2, 0, 3, 0),
  @(5002, "soul.png", 1, "sk_human,charactor2/face.png,charactor2/head.png,charactor2/body.png,charactor2/arm.png,charactor2/hand.png,charactor2/arm.png,charactor2/hand.png,charactor2/leg.png,charactor15/leg.png,", 0, "SOUL_NAME_5002", "SOUL_DESC_5002", 1, 5, 1, 0, 3000, 50, 40, 10, 10, 10, 0, 0, 1, 0, 4, 0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $rowVals = $data[$i]
  for ($c = 1; $c -le $rowVals.Count; $c++) {
    $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
  }

  # Row 2 (first soul) keeps its original compact formatting; every other
  # data row (3 and on) gets the taller wrap-text row used for the long
  # avatar strings in column D.
  if ($r -gt 2) {
    $ws.Range("D" + $r).WrapText = $true
    $ws.Rows.Item($r).RowHeight = 27
  }
}

# ---------------------------------------------------------------------
# 3. Cosmetic sheet-view tweaks: widen column C, move the selection.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 14.0357142857
$ws.Range("O8").Select()
